# This workbook is a "weekly" price log for "Espinaca" (spinach) at the
# "Vega Central Mapocho de Santiago" market. Each week two new rows (one
# for "Primera" quality, one for "Segunda" quality) are prepended to the
# top of the data block, pushing the older rows down.
#
# Data block in this sheet runs from row 2 (first data row, header is
# row 1) through row 229 (before the edit). The two new rows are
# inserted right before what is currently row 153, i.e. they become the
# new rows 153-154, and the previous rows 153-229 shift down to 155-231.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new blank rows above the current row 153 -----------
$ws.Rows.Item(153).EntireRow.Insert()
$ws.Rows.Item(153).EntireRow.Insert()

# --- 2. Fill in the new row 153 ("Primera") ----------------------------
$ws.Cells.Item(153, 1).Value  = 9
$ws.Cells.Item(153, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(153, 3).Value  = "Metropolitana"
$ws.Cells.Item(153, 4).Value  = 44460
$ws.Cells.Item(153, 5).Value  = 13
$ws.Cells.Item(153, 6).Value  = 100112012
$ws.Cells.Item(153, 7).Value  = "Espinaca"
$ws.Cells.Item(153, 8).Value  = "Sin especificar"
$ws.Cells.Item(153, 9).Value  = "Primera"
$ws.Cells.Item(153, 10).Value = 196
$ws.Cells.Item(153, 11).Value = 6000
$ws.Cells.Item(153, 12).Value = 7000
$ws.Cells.Item(153, 13).Value = 6500
$ws.Cells.Item(153, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(153, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(153, 16).Value = 650
$ws.Cells.Item(153, 17).Value = 10
$ws.Cells.Item(153, 18).Value = "Hortaliza"
$ws.Cells.Item(153, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- 3. Fill in the new row 154 ("Segunda") ----------------------------
$ws.Cells.Item(154, 1).Value  = 9
$ws.Cells.Item(154, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(154, 3).Value  = "Metropolitana"
$ws.Cells.Item(154, 4).Value  = 44460
$ws.Cells.Item(154, 5).Value  = 13
$ws.Cells.Item(154, 6).Value  = 100112012
$ws.Cells.Item(154, 7).Value  = "Espinaca"
$ws.Cells.Item(154, 8).Value  = "Sin especificar"
$ws.Cells.Item(154, 9).Value  = "Segunda"
$ws.Cells.Item(154, 10).Value = 97
$ws.Cells.Item(154, 11).Value = 5000
$ws.Cells.Item(154, 12).Value = 5000
$ws.Cells.Item(154, 13).Value = 5000
$ws.Cells.Item(154, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(154, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(154, 16).Value = 500
$ws.Cells.Item(154, 17).Value = 10
$ws.Cells.Item(154, 18).Value = "Hortaliza"
$ws.Cells.Item(154, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
